$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2022-11-11 10:18:01.798523"
$ws.Range("B2").Value = "Calle 35, Santa Teresita, Comuna 12 - La América, Perímetro Urbano Medellín, Medellín, Valle de Aburrá, Antioquia, 050032, Colombia"
